$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 3850.2,               4202, 3414, 0.08217678864796957),
    @(1, 3459.533333333333,    3786, 2983, 0.08249160448710123),
    @(2, 3923.9,                4223, 3592, 0.08638970851898194),
    @(3, 3642.4,                3939, 3327, 0.08260746796925862),
    @(4, 3577.066666666667,    3903, 3147, 0.08259633382161459),
    @(5, 3757.633333333333,    4088, 3293, 0.08684235413869222),
    @(6, 4125.133333333333,    4397, 3694, 0.08592217763264974),
    @(7, 3563,                  3869, 3356, 0.08363552093505859),
    @(8, 3881.866666666667,    4150, 3362, 0.08337690830230712),
    @(9, 3789.766666666667,    4217, 3319, 0.0830965518951416)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
